# Regenerate s_vals data to filter save games.
# Updates the computed statistic columns (TB, d2S, K, IP, sum) for each
# game row on Sheet1, leaving the date (A) and Win (F) columns untouched
# in value (F values happen to stay the same, but are rewritten for
# consistency with the source regeneration).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(0.1554434735375247, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1, 1.705647867635037)
    3 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
    4 = @(0.3464964993005633, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 0, 1.611428759079651)
    5 = @(3.182878228561681, 0.3375848360084654, 0.7127328510149897, 0.4998867070740569, 1, 4.733082622659194)
    6 = @(1.505614041169197, 1766.335244827366, 0.7127328510149897, 6.48142807727062, 1, 1775.035019796821)
    7 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 1, 8.418600821238126)
    8 = @(3.182878228561681, 1.65323645889881, 16.98373111632243, 0.4998867070740569, 1, 22.31973251085698)
}

$columns = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
